# Updates cryptos list values (Price and Volume(1h) columns) as per
# "Updated cryptos list on Sat Sep 28 09:40:43 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '65.568.30'
$ws.Range('E2').Value = '  -0.29%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.666.01'
$ws.Range('E3').Value = '  +0.23%  '

# Row 4
$ws.Range('E4').Value = '  +0.00%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '599.86'
$ws.Range('E5').Value = '  -1.16%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '156.44'
$ws.Range('E6').Value = '  -1.06%  '

# Row 7
$ws.Range('E7').Value = '  +0.03%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.604'
$ws.Range('E8').Value = '  +2.49%  '

# Row 9
$ws.Range('E9').Value = '  -1.85%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.91'
$ws.Range('E10').Value = '  +0.62%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.397'
$ws.Range('E11').Value = '  -1.95%  '

# Row 12
$ws.Range('E12').Value = '  -0.17%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '29.35'
$ws.Range('E13').Value = '  -1.97%  '

# Row 14
$ws.Range('E14').Value = '  -0.42%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.145.57'
$ws.Range('E15').Value = '  +0.14%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.394.84'
$ws.Range('E16').Value = '  -0.11%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.657.63'
$ws.Range('E17').Value = '  +0.43%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.48'
$ws.Range('E18').Value = '  -1.90%  '

# Row 19
$ws.Range('E19').Value = '  -1.74%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.49'
$ws.Range('E20').Value = '  +0.76%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '350.33'
$ws.Range('E21').Value = '  -3.17%  '

# Row 22
$ws.Range('E22').Value = '  -0.12%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '69.68'
$ws.Range('E23').Value = '  +0.35%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.66'
$ws.Range('E24').Value = '  +0.70%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000108'
$ws.Range('E25').Value = '  +1.87%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.64'
$ws.Range('E26').Value = '  -3.69%  '

# Row 27
$ws.Range('E27').Value = '  +0.87%  '

# Row 28
$ws.Range('E28').Value = '  -3.44%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.07'
$ws.Range('E29').Value = '  -1.43%  '

# Row 30
$ws.Range('B30').Value = 'Binance-PegBSC-USD'
$ws.Range('C30').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.00'
$ws.Range('E30').Value = '  -0.08%  '

# Row 31
$ws.Range('B31').Value = 'Bittensor'
$ws.Range('C31').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '539.12'
$ws.Range('E31').Value = '  -0.20%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.14'
$ws.Range('E32').Value = '  -3.27%  '

# Row 33
$ws.Range('E33').Value = '  -4.67%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.53'
$ws.Range('E34').Value = '  +2.46%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.44'
$ws.Range('E35').Value = '  -3.01%  '

# Row 36
$ws.Range('E36').Value = '  -3.01%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '20.41'
$ws.Range('E37').Value = '  -1.37%  '

# Row 38
$ws.Range('E38').Value = '  -0.05%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '158.79'
$ws.Range('E39').Value = '  -2.41%  '

# Row 40
$ws.Range('E40').Value = '  -3.62%  '

# Row 41
$ws.Range('E41').Value = '  +0.02%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '42.47'
$ws.Range('E42').Value = '  +0.17%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '165.46'
$ws.Range('E43').Value = '  -0.39%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.07'
$ws.Range('E44').Value = '  -2.55%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0610'
$ws.Range('E45').Value = '  -0.72%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.26'
$ws.Range('E46').Value = '  -5.68%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '22.99'
$ws.Range('E47').Value = '  -0.50%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.646'
$ws.Range('E48').Value = '  -2.20%  '

# Row 49
$ws.Range('E49').Value = '  -1.99%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0996'
$ws.Range('E50').Value = '  +0.68%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '19.99'
$ws.Range('E51').Value = '  +0.86%  '
